$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 10.02.2022 00:15"

# Row 10 ("EuroOil Opuštěná") price update:
#  - New current price (B10) becomes 36.9
#  - Old current price moves into "Old Cena" (C10) as 36.6
#  - Delta (D10) becomes the text "+0.3" instead of a plain number
#  - Old Datum (E10) becomes a plain text timestamp instead of a formatted date serial
$ws.Range("B10").Value = 36.9
$ws.Range("C10").Value = 36.6

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "+0.3"
$ws.Range("D10").ClearFormats()

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2022-02-10 00:20:52"
$ws.Range("E10").ClearFormats()
